$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in attendance for the week of 2/24 (column F) for each person.
$ws.Range("F3").Value = 1   # Aashay
$ws.Range("F4").Value = 1   # Saurabh
$ws.Range("F5").Value = 1   # Ania
$ws.Range("F6").Value = 1   # Attush
$ws.Range("F7").Value = 0   # Cindy
$ws.Range("F8").Value = 1   # Lindsey

# Move the active selection to D12 (matches the author's saved cursor position).
$ws.Range("D12").Select()
